$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain plain text so Excel does not
# auto-convert numeric-looking strings (e.g. "1.00") into numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "37.766.83"
$ws.Range("E2").Value = "  +2.11%  "
$ws.Range("D3").Value = "2.045.77"
$ws.Range("E3").Value = "  +1.70%  "
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").Value = "230.24"
$ws.Range("E5").Value = "  +1.85%  "
$ws.Range("D6").Value = "0.614"
$ws.Range("E6").Value = "  +1.07%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "56.55"
$ws.Range("E8").Value = "  +3.56%  "
$ws.Range("D9").Value = "0.384"
$ws.Range("E9").Value = "  +1.39%  "
$ws.Range("E10").Value = "  +2.71%  "
$ws.Range("E11").Value = "  -0.91%  "
$ws.Range("D12").Value = "2.346.84"
$ws.Range("E12").Value = "  +1.68%  "
$ws.Range("D13").Value = "14.47"
$ws.Range("E13").Value = "  +1.47%  "
$ws.Range("D14").Value = "20.41"
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").Value = "5.24"
$ws.Range("E15").Value = "  +2.35%  "
$ws.Range("D16").Value = "0.744"
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("D17").Value = "2.033.14"
$ws.Range("E17").Value = "  +1.13%  "
$ws.Range("D18").Value = "37.668.78"
$ws.Range("E18").Value = "  +2.16%  "
$ws.Range("D19").Value = "6.20"
$ws.Range("E19").Value = "  +0.94%  "
$ws.Range("D20").Value = "69.32"
$ws.Range("E20").Value = "  +0.84%  "
$ws.Range("D21").Value = "0.0₃0826"
$ws.Range("E21").Value = "  +0.92%  "
$ws.Range("D22").Value = "224.23"
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  +1.43%  "
$ws.Range("E25").Value = "  +3.77%  "
$ws.Range("D26").Value = "165.48"
$ws.Range("E26").Value = "  +0.44%  "
$ws.Range("D27").Value = "9.23"
$ws.Range("E27").Value = "  +0.39%  "
$ws.Range("D28").Value = "0.133"
$ws.Range("E28").Value = "  +5.80%  "
$ws.Range("D29").Value = "19.08"
$ws.Range("E29").Value = "  +2.09%  "
$ws.Range("D30").Value = "1.34"
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("E31").Value = "  +0.91%  "
$ws.Range("D32").Value = "4.51"
$ws.Range("E32").Value = "  +0.31%  "
$ws.Range("D33").Value = "0.0611"
$ws.Range("E33").Value = "  -0.35%  "
$ws.Range("D34").Value = "4.52"
$ws.Range("E34").Value = "  +2.61%  "
$ws.Range("E35").Value = "  +9.27%  "
$ws.Range("E36").Value = "  +0.37%  "
$ws.Range("D37").Value = "5.93"
$ws.Range("E37").Value = "  +11.12%  "
$ws.Range("D38").Value = "3.27"
$ws.Range("E38").Value = "  +4.13%  "
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("D40").Value = "98.27"
$ws.Range("E40").Value = "  +3.26%  "
$ws.Range("D41").Value = "1.482.01"
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("D42").Value = "0.0951"
$ws.Range("E42").Value = "  +3.35%  "
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("E44").Value = "  +2.54%  "
$ws.Range("D45").Value = "16.66"
$ws.Range("E45").Value = "  -0.69%  "
$ws.Range("D46").Value = "4.14"
$ws.Range("E46").Value = "  +14.84%  "
$ws.Range("D47").Value = "1.12"
$ws.Range("E47").Value = "  -0.79%  "
$ws.Range("D48").Value = "1.02"
$ws.Range("E48").Value = "  +1.28%  "
$ws.Range("E49").Value = "  -1.03%  "
$ws.Range("E50").Value = "  +0.56%  "
$ws.Range("D51").Value = "2.231.39"
$ws.Range("E51").Value = "  +1.57%  "

# Restore the default style on column D now that values are set,
# so no stray explicit cell format lingers on these cells.
$priceRange.Style = "Normal"
